# Coefficients workbook update:
#  - add 3 new rows (14, 15, 16) of cluster/state coefficients ("0.5" and "0.01")
#  - update the active selection to C17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New numeric ids for the new rows
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(16, 1).Value = 15

# Values for the new rows, column by column (B..I == col 2..9)
# Row 14: all columns = "0.5"
# Row 15: B..E = "0.5", F..I = "0.01"
# Row 16: B..E = "0.01", F..I = "0.5"
$rowValues = @{
    14 = @("0.5", "0.5", "0.5", "0.5", "0.5", "0.5", "0.5", "0.5")
    15 = @("0.5", "0.5", "0.5", "0.5", "0.01", "0.01", "0.01", "0.01")
    16 = @("0.01", "0.01", "0.01", "0.01", "0.5", "0.5", "0.5", "0.5")
}

# Force the cells to hold real text values (so "0.5"/"0.01" are written as
# shared strings, matching the rest of the coefficient table) rather than
# being auto-converted to numbers, then strip the temporary text format
# back off so no stray cell styling is left behind.
$targetRange = $ws.Range("B14:I16")
$targetRange.NumberFormat = "@"

foreach ($r in 14..16) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

$targetRange.ClearFormats()

# Update the selected cell shown when the sheet is opened
$ws.Range("C17").Select() | Out-Null
